$wb = $excel.ActiveWorkbook

# Update the raw metric values on the "Metrics" sheet (B2:B13).
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 325202.15000000002
$metrics.Range("B3").Value  = 278595.95999999996
$metrics.Range("B4").Value  = 99377.459999999992
$metrics.Range("B5").Value  = 13280
$metrics.Range("B6").Value  = 5527909.2600000007
$metrics.Range("B7").Value  = 4678948.9200000009
$metrics.Range("B8").Value  = 1631334.34
$metrics.Range("B9").Value  = 215987
$metrics.Range("B10").Value = 33993290.249999993
$metrics.Range("B11").Value = 31954224.079999998
$metrics.Range("B12").Value = 11913056.379999995
$metrics.Range("B13").Value = 1313617

# Move the active selection on the "today" sheet to F6, matching the
# author's recorded cursor position at save time.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F6").Select()
